$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price and Volume columns remain plain text so values like "1.00" or "0.999" are preserved exactly
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range('D2').Value = '63.589.89'
$ws.Range('E2').Value = '  -6.53%  '
$ws.Range('D3').Value = '3.258.07'
$ws.Range('E3').Value = '  -8.90%  '
$ws.Range('D4').Value = '1.00'
$ws.Range('E4').Value = '  -0.09%  '
$ws.Range('D5').Value = '175.20'
$ws.Range('E5').Value = '  -14.56%  '
$ws.Range('D6').Value = '510.61'
$ws.Range('E6').Value = '  -9.58%  '
$ws.Range('D7').Value = '0.585'
$ws.Range('E7').Value = '  -4.26%  '
$ws.Range('B8').Value = 'USDC'
$ws.Range('C8').Value = 'https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc'
$ws.Range('D8').Value = '1.00'
$ws.Range('E8').Value = '  -0.03%  '
$ws.Range('B9').Value = 'LidoStakedEther'
$ws.Range('C9').Value = 'https://coinranking.com/coin/VINVMYf0u+lidostakedether-steth'
$ws.Range('D9').Value = '3.253.84'
$ws.Range('E9').Value = '  -8.93%  '
$ws.Range('D10').Value = '0.606'
$ws.Range('E10').Value = '  -10.43%  '
$ws.Range('D11').Value = '56.50'
$ws.Range('E11').Value = '  -6.20%  '
$ws.Range('D12').Value = '0.129'
$ws.Range('E12').Value = '  -12.00%  '
$ws.Range('D13').Value = '0.0000251'
$ws.Range('E13').Value = '  -9.79%  '
$ws.Range('D14').Value = '8.91'
$ws.Range('E14').Value = '  -12.53%  '
$ws.Range('D15').Value = '3.801.26'
$ws.Range('E15').Value = '  -8.55%  '
$ws.Range('D16').Value = '0.118'
$ws.Range('E16').Value = '  -6.85%  '
$ws.Range('D17').Value = '3.269.73'
$ws.Range('E17').Value = '  -9.07%  '
$ws.Range('D18').Value = '63.368.94'
$ws.Range('E18').Value = '  -6.65%  '
$ws.Range('D19').Value = '16.87'
$ws.Range('E19').Value = '  -10.27%  '
$ws.Range('D20').Value = '10.65'
$ws.Range('E20').Value = '  -12.42%  '
$ws.Range('D21').Value = '0.934'
$ws.Range('E21').Value = '  -11.56%  '
$ws.Range('D22').Value = '364.05'
$ws.Range('E22').Value = '  -9.04%  '
$ws.Range('D23').Value = '79.09'
$ws.Range('E23').Value = '  -5.99%  '
$ws.Range('D24').Value = '3.59'
$ws.Range('E24').Value = '  -13.70%  '
$ws.Range('D25').Value = '10.68'
$ws.Range('E25').Value = '  -14.06%  '
$ws.Range('D26').Value = '3.76'
$ws.Range('E26').Value = '  -2.92%  '
$ws.Range('D27').Value = '5.97'
$ws.Range('E27').Value = '  -2.55%  '
$ws.Range('D28').Value = '2.59'
$ws.Range('E28').Value = '  -9.42%  '
$ws.Range('D29').Value = '11.12'
$ws.Range('E29').Value = '  -10.12%  '
$ws.Range('D30').Value = '8.17'
$ws.Range('E30').Value = '  -10.85%  '
$ws.Range('D31').Value = '28.04'
$ws.Range('E31').Value = '  -10.54%  '
$ws.Range('D32').Value = '628.98'
$ws.Range('E32').Value = '  -4.68%  '
$ws.Range('D33').Value = '6.53'
$ws.Range('E33').Value = '  -15.08%  '
$ws.Range('D34').Value = '10.92'
$ws.Range('E34').Value = '  -9.16%  '
$ws.Range('D35').Value = '58.99'
$ws.Range('E35').Value = '  -6.69%  '
$ws.Range('E36').Value = '  -9.63%  '
$ws.Range('E37').Value = '  +0.07%  '
$ws.Range('D38').Value = '35.12'
$ws.Range('E38').Value = '  -14.24%  '
$ws.Range('D39').Value = '0.368'
$ws.Range('E39').Value = '  -9.62%  '
$ws.Range('D40').Value = '1.00'
$ws.Range('E40').Value = '  +0.20%  '
$ws.Range('B41').Value = 'Maker'
$ws.Range('C41').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D41').Value = '2.828.73'
$ws.Range('E41').Value = '  -10.91%  '
$ws.Range('B42').Value = 'Kaspa'
$ws.Range('C42').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D42').Value = '0.120'
$ws.Range('E42').Value = '  -9.57%  '
$ws.Range('D43').Value = '0.0₃0640'
$ws.Range('E43').Value = '  -14.28%  '
$ws.Range('D44').Value = '2.63'
$ws.Range('E44').Value = '  -19.17%  '
$ws.Range('E45').Value = '  -8.09%  '
$ws.Range('D46').Value = '0.0376'
$ws.Range('E46').Value = '  -8.10%  '
$ws.Range('D47').Value = '2.26'
$ws.Range('E47').Value = '  -15.88%  '
$ws.Range('D48').Value = '0.122'
$ws.Range('E48').Value = '  -6.12%  '
$ws.Range('D49').Value = '131.74'
$ws.Range('E49').Value = '  -4.96%  '
$ws.Range('B50').Value = 'ApeXProtocol'
$ws.Range('C50').Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range('D50').Value = '2.84'
$ws.Range('E50').Value = '  -6.93%  '
$ws.Range('B51').Value = 'Stacks'
$ws.Range('C51').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D51').Value = '2.62'
$ws.Range('E51').Value = '  -2.26%  '
